$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$avgRange = $ws.Range("B14:B17")
$avgRange.Font.Bold = $true
$avgRange.Font.Size = 12
$avgRange.VerticalAlignment = -4108

$ws.Range("J12").Font.Bold = $true
